$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row total for Right column (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update "Total" row for Right column (B12: 45 -> 75)
$ws.Range("B12").Value = 75

# Update Correct/Total marks text (E12: "42/84" -> "75/140")
$ws.Range("E12").Value = "75/140"
